# Applies the commit: add a new "Player Info" sheet (as the first sheet)
# and update the existing "ODI Batting" sheet's MATCH_CARD_LINK column to a
# plain MATCH_CODE value.
#
# NOTE: worksheet object references here resolve by live position, so all
# work on the pre-existing "ODI Batting" sheet is done BEFORE any new sheet
# is inserted (inserting a sheet shifts indices out from under a previously
# captured reference).

$wb = $excel.ActiveWorkbook
$batting = $wb.ActiveSheet

# --- 1. Update the "ODI Batting" sheet ---------------------------------------
# Rename the MATCH_CARD_LINK header to MATCH_CODE.
$batting.Range("D1").Value = "MATCH_CODE"

# Replace the full scorecard URLs with the bare numeric match code, keeping
# the value stored as text (as the rest of the sheet's cells are) and
# stripping the temporary number-format back off afterwards.
$codeRange = $batting.Range("D2:D4")
$codeRange.NumberFormat = "@"
$batting.Range("D2").Value = "4519"
$batting.Range("D3").Value = "4520"
$batting.Range("D4").Value = "4522"
$codeRange.Style = "Normal"

# --- 2. Insert a brand-new worksheet in front of it and name it -------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Match the page margins used throughout the rest of the workbook (points:
# 0.75in/1in/0.5in == 54/72/36pt) instead of the engine's bare defaults.
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

# --- 3. Populate the new "Player Info" sheet ---------------------------------
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row - force text storage (IDs look numeric but must stay text), then
# drop back to the Normal style so no stray number-format sticks around.
$playerInfoData = $playerInfo.Range("A2:D2")
$playerInfoData.NumberFormat = "@"
$playerInfo.Range("A2").Value = "6148"
$playerInfo.Range("B2").Value = "Neil Alan Rock"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"
$playerInfoData.Style = "Normal"

# Header formatting to match the bold / bordered / centered header style
# used on row 1 of the "ODI Batting" sheet.
$playerInfoHeader = $playerInfo.Range("A1:D1")
$playerInfoHeader.Font.Bold = $true
$playerInfoHeader.HorizontalAlignment = -4108
$playerInfoHeader.VerticalAlignment = -4160
$playerInfoHeader.Borders.Item(1).LineStyle = 1
$playerInfoHeader.Borders.Item(1).Weight = 2
$playerInfoHeader.Borders.Item(2).LineStyle = 1
$playerInfoHeader.Borders.Item(2).Weight = 2
$playerInfoHeader.Borders.Item(3).LineStyle = 1
$playerInfoHeader.Borders.Item(3).Weight = 2
$playerInfoHeader.Borders.Item(4).LineStyle = 1
$playerInfoHeader.Borders.Item(4).Weight = 2

Write-Host "Workbook updated: added 'Player Info' sheet and refreshed MATCH_CODE column."
